$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (week number + date range) ---
$ws.Range("A8").Characters(21, 2).Text = "50"
$ws.Range("C9").Characters(47, 9).Text = "12/14/2025"
$ws.Range("C9").Characters(27, 9).Text = "12/8/2025"

# --- Simple numeric value updates (style/type unchanged) ---
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("I15").Value = 14
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = 7.692307692307
$ws.Range("L15").Value = 133.333333333333
$ws.Range("M15").Value = 1300
$ws.Range("N15").Value = 75
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 8
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 111
$ws.Range("J16").Value = 90
$ws.Range("K16").Value = 23.333333333333
$ws.Range("L16").Value = -5.932203389830
$ws.Range("M16").Value = -10.483870967741
$ws.Range("N16").Value = -80.862068965517
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -27.272727272727
$ws.Range("I17").Value = 138
$ws.Range("J17").Value = 132
$ws.Range("K17").Value = 4.545454545454
$ws.Range("L17").Value = 27.777777777777
$ws.Range("M17").Value = 94.366197183098
$ws.Range("N17").Value = -45.669291338582
$ws.Range("F18").Value = 3
$ws.Range("H18").Value = -57.142857142857
$ws.Range("I18").Value = 129
$ws.Range("J18").Value = 141
$ws.Range("K18").Value = -8.510638297872
$ws.Range("L18").Value = -23.214285714285
$ws.Range("M18").Value = -44.635193133047
$ws.Range("N18").Value = -88.088642659279
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 11.111111111111
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = -16.981132075471
$ws.Range("I19").Value = 578
$ws.Range("J19").Value = 654
$ws.Range("K19").Value = -11.620795107033
$ws.Range("L19").Value = -11.485451761102
$ws.Range("M19").Value = 97.269624573378
$ws.Range("N19").Value = 64.204545454545
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 119
$ws.Range("J20").Value = 109
$ws.Range("K20").Value = 9.174311926605
$ws.Range("L20").Value = -20.666666666666
$ws.Range("M20").Value = -15.602836879432
$ws.Range("N20").Value = -86.353211009174
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -15.789473684210
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = -13.253012048192
$ws.Range("I21").Value = 1089
$ws.Range("J21").Value = 1139
$ws.Range("K21").Value = -4.389815627743
$ws.Range("L21").Value = -9.476309226932
$ws.Range("M21").Value = 26.187717265353
$ws.Range("N21").Value = -65.505226480836
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 34
$ws.Range("J23").Value = 29
$ws.Range("K23").Value = 17.241379310344
$ws.Range("L23").Value = 36
$ws.Range("M23").Value = 54.545454545454
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -13.333333333333
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = -24.324324324324
$ws.Range("I24").Value = 896
$ws.Range("J24").Value = 947
$ws.Range("K24").Value = -5.385427666314
$ws.Range("L24").Value = 4.795321637426
$ws.Range("M24").Value = 53.162393162393
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = -43.478260869565
$ws.Range("I25").Value = 431
$ws.Range("J25").Value = 597
$ws.Range("K25").Value = -27.805695142378
$ws.Range("L25").Value = -8.102345415778
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 5.882352941176
$ws.Range("I26").Value = 266
$ws.Range("J26").Value = 263
$ws.Range("K26").Value = 1.140684410646
$ws.Range("L26").Value = 13.675213675213
$ws.Range("M26").Value = 31.683168316831
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("I27").Value = 17
$ws.Range("J27").Value = 14
$ws.Range("K27").Value = 21.428571428571
$ws.Range("L27").Value = 112.5
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 47
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 17.5
$ws.Range("L28").Value = 9.302325581395

# --- Cells whose style/type flips between numeric and shared-text ---
# D22: numeric 1 -> text "0" (style 13, like C22)
$ws.Range("D22").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("D22").PasteSpecial(-4122)

# E22: numeric -100 -> text "***.*" (style 13, like N22)
$ws.Range("E22").Value = "'***.*"
$ws.Range("N22").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# C23: text "0" -> numeric 1 (style 15, like D23)
$ws.Range("C23").Value = 1
$ws.Range("D23").Copy()
$ws.Range("C23").PasteSpecial(-4122)

# D28: text "0" -> numeric 1 (style 15, like C28)
$ws.Range("D28").Value = 1
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)

# E28: text "***.*" -> numeric 0 (style 14, like H28)
$ws.Range("E28").Value = 0
$ws.Range("H28").Copy()
$ws.Range("E28").PasteSpecial(-4122)

# F31: numeric 1 -> text "0" (style 13, like G31)
$ws.Range("F31").Value = "'0"
$ws.Range("G31").Copy()
$ws.Range("F31").PasteSpecial(-4122)

